$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S3").Value = 1.69
$ws.Range("S4").Value = 1.69
$ws.Range("G5").Value = 1.65
$ws.Range("N5").Value = 1.85
$ws.Range("O5").Value = 2
$ws.Range("R5").Value = 1.8
$ws.Range("S5").Value = 1.8
$ws.Range("R7").Value = 1.81
$ws.Range("S7").Value = 1.9
$ws.Range("N8").Value = 1.58
$ws.Range("G9").Value = 1.77
$ws.Range("I9").Value = 3.7
$ws.Range("M9").Value = 3.8
$ws.Range("N9").Value = 1.7
$ws.Range("O9").Value = 2.1
$ws.Range("P9").Value = 1.35
$ws.Range("Q9").Value = 3
$ws.Range("S9").Value = 2.1
$ws.Range("T9").Value = 8.5
$ws.Range("U9").Value = 9.5
$ws.Range("W9").Value = 15
$ws.Range("P10").Value = 1.34
$ws.Range("Q10").Value = 3.1
$ws.Range("Q12").Value = 2.77
$ws.Range("G13").Value = 2.02
$ws.Range("H13").Value = 3.15
$ws.Range("I13").Value = 3.5
$ws.Range("K13").Value = 6.4
$ws.Range("L13").Value = 1.39
$ws.Range("M13").Value = 2.77
$ws.Range("N13").Value = 2.15
$ws.Range("O13").Value = 1.62
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.42
$ws.Range("U13").Value = 9
$ws.Range("W13").Value = 18
$ws.Range("X13").Value = 18
$ws.Range("Z13").Value = 6.4
$ws.Range("AA13").Value = 6.2
$ws.Range("AB13").Value = 16
$ws.Range("AE13").Value = 9.25
$ws.Range("AF13").Value = 18.5
$ws.Range("AH13").Value = 50
$ws.Range("AI13").Value = 35
$ws.Range("P14").Value = 1.3
$ws.Range("R14").Value = 1.67
$ws.Range("G17").Value = 2.25
$ws.Range("R17").Value = 1.72
$ws.Range("G18").Value = 1.25
$ws.Range("R18").Value = 1.8
$ws.Range("S18").Value = 1.8
$ws.Range("G19").Value = 1.77
$ws.Range("L19").Value = 1.24
$ws.Range("Q19").Value = 2.9
$ws.Range("S19").Value = 2.12
$ws.Range("R20").Value = 1.69
$ws.Range("G22").Value = 1.53
$ws.Range("H22").Value = 3.95
$ws.Range("I22").Value = 5.7
$ws.Range("K22").Value = 7.5
$ws.Range("L22").Value = 1.29
$ws.Range("M22").Value = 3.3
$ws.Range("N22").Value = 1.87
$ws.Range("O22").Value = 1.85
$ws.Range("Q22").Value = 2.67
$ws.Range("R22").Value = 1.95
$ws.Range("S22").Value = 1.75
$ws.Range("U22").Value = 6.8
$ws.Range("W22").Value = 10.5
$ws.Range("Z22").Value = 7.5
$ws.Range("AA22").Value = 7.6
$ws.Range("AB22").Value = 18.5
$ws.Range("AE22").Value = 14
$ws.Range("AF22").Value = 32
$ws.Range("AG22").Value = 18
$ws.Range("AH22").Value = 110
$ws.Range("AI22").Value = 60
$ws.Range("AJ22").Value = 65
